$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments ---
# (Excel COM's ColumnWidth is in character units and the engine adds a
#  fixed 5/6 character padding when serializing the raw OOXML <col> width,
#  so subtract that padding here to land on the exact target widths.)
$ws.Columns.Item(3).ColumnWidth = 71 - 5/6
$ws.Columns.Item(4).ColumnWidth = 44 - 5/6
$ws.Columns.Item(6).ColumnWidth = 17 - 5/6

# --- Data rows (table body) ---
$data = @(
    @("1324475", "https://aiesec.org/opportunity/global-talent/1324475", "Customer Experience Intern", "Maastricht, Netherlands", "Yes", "187 applicants", "6 - 18 Months", "DHL Group"),
    @("1326517", "https://aiesec.org/opportunity/global-talent/1326517", "Researcher", "Panamá, Provincia de Panamá, Panamá", "No", "0 applicants", "6 - 18 Months", "Michael Page International Panamá S.A."),
    @("1326505", "https://aiesec.org/opportunity/global-talent/1326505", "Accelerate Romania - Digital Content Intern – Social Media & Website", "Cluj-Napoca, Romania", "No", "4 applicants", "9 - 12 Weeks", "Dog Assist"),
    @("1326504", "https://aiesec.org/opportunity/global-talent/1326504", "BI RPM Intern", "Panamá, Provincia de Panamá, Panamá", "No", "4 applicants", "6 - 18 Months", "Samsung Electronics Latinoamérica (Zona Libre) S.A (SELA)"),
    @("1326501", "https://aiesec.org/opportunity/global-talent/1326501", "People Data Specialist Intern", "Fritz-Erler-Straße 5, 53113 Bonn, Germany", "Yes", "10 applicants", "6 - 18 Months", "DHL Group"),
    @("1326481", "https://aiesec.org/opportunity/global-talent/1326481", "Global Duty Billing Data Analytics Expert", "Maastricht, Netherlands", "Yes", "14 applicants", "6 - 18 Months", "DHL Group"),
    @("1289255", "https://aiesec.org/opportunity/global-talent/1289255", "Medical Advisor French Speaker", "Fatih, Türkiye", "No", "644 applicants", "6 - 18 Months", "International Plus")
)

$rowIndex = 2
foreach ($row in $data) {
    for ($col = 1; $col -le 8; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $row[$col - 1]
    }
    $rowIndex++
}

# --- Highlight the PREMIUM = Yes cells with a yellow fill ---
$ws.Cells.Item(2, 5).Interior.Color = 65535
$ws.Cells.Item(6, 5).Interior.Color = 65535
$ws.Cells.Item(7, 5).Interior.Color = 65535

